$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly observation was inserted at row 48 (Fecha 44677), pushing every
# following record down by one row (old row 48 -> new row 49, ..., old row 306
# -> new row 307). Insert a whole row so Excel shifts everything down for us.
$ws.Rows.Item(48).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A48").Value = 3
$ws.Range("B48").Value = "Femacal de La Calera"
$ws.Range("C48").Value = "Coquimbo"
$ws.Range("D48").Value = 44677
$ws.Range("E48").Value = 5
$ws.Range("F48").Value = 100112039
$ws.Range("G48").Value = "Ciboulette"
$ws.Range("H48").Value = "Sin especificar"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 120
$ws.Range("K48").Value = 1500
$ws.Range("L48").Value = 1500
$ws.Range("M48").Value = 1500
$ws.Range("N48").Value = "`$/docena de atados"
$ws.Range("O48").Value = "Provincia de Quillota"
$ws.Range("P48").Value = 500
$ws.Range("Q48").Value = 3
$ws.Range("R48").Value = "Hortaliza"
